$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddProductCategory1")

# --- Row 1: remove J1 (20), keep K1 (30) ---
$ws.Range("J1").ClearContents()

# --- Row 2: overwrite G2 with the 1970-01-01 date (same as G1/G3), remove J2 ---
$ws.Range("G2").Value = 25569
$ws.Range("G2").NumberFormat = "DD/MM/YY"
$ws.Range("J2").ClearContents()

# --- Row 3: add G3 with the same 1970-01-01 date, remove J3 ---
$ws.Range("G3").Value = 25569
$ws.Range("G3").NumberFormat = "DD/MM/YY"
$ws.Range("J3").ClearContents()

# --- Row 4: new "mid-2018" currency snapshot rows (G/H/I only) ---
$ws.Range("G4").Value = 43276.5385855949
$ws.Range("G4").NumberFormat = "DD/MM/YY"
$ws.Range("H4").NumberFormat = "DD/MM/YY"
$ws.Range("H4").Value = "United States Dollar"
$ws.Range("I4").Value = 20

$ws.Range("G5").Value = 43276.5385094428
$ws.Range("G5").NumberFormat = "DD/MM/YY"
$ws.Range("H5").Value = "Canadian Dollar"
$ws.Range("I5").Value = 21

$ws.Range("G6").Value = 43276.5385244457
$ws.Range("G6").NumberFormat = "DD/MM/YY"
$ws.Range("H6").Value = "Euro"
$ws.Range("I6").Value = 22

# --- Rows 7-9: new "2020-01-01" currency snapshot rows (G/H/I only) ---
$ws.Range("G7").Value = 43831
$ws.Range("G7").NumberFormat = "DD/MM/YY"
$ws.Range("H7").NumberFormat = "DD/MM/YY"
$ws.Range("H7").Value = "United States Dollar"
$ws.Range("I7").Value = 30

$ws.Range("G8").Value = 43831
$ws.Range("G8").NumberFormat = "DD/MM/YY"
$ws.Range("H8").Value = "Canadian Dollar"
$ws.Range("I8").Value = 31

$ws.Range("G9").Value = 43831
$ws.Range("G9").NumberFormat = "DD/MM/YY"
$ws.Range("H9").Value = "Euro"
$ws.Range("I9").Value = 32

# --- Selection moves to I19 ---
$ws.Range("I19").Select()

Write-Output "done"
